{"js": "// Add a new bulleted list item \"Twig\" right after the \"Npm 9.5.1\" item\n// in the \"Pour le front :\" list (same list/style as its neighbours).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text,items/style\");\nawait context.sync();\n\n// Locate the \"Npm 9.5.1\" list item (unique text fragment \"9.5.1\").\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  if (para.text && para.text.indexOf(\"9.5.1\") !== -1) {\n    target = para;\n    break;\n  }\n}\nif (!target) {\n  throw new Error('Could not find the \"9.5.1\" paragraph to anchor the new item.');\n}\n\n// Read the list this paragraph belongs to, so the new paragraph can join it.\nconst list = target.list;\nlist.load(\"id\");\nawait context.sync();\n\n// Insert \"Twig\" as a new paragraph right after it, matching style + list.\nconst newPara = target.insertParagraph(\"Twig\", Word.InsertLocation.after);\nnewPara.style = target.style;\nnewPara.attachToList(list.id, 0);\n\nawait context.sync();\n", "ps1": "# Add a new bulleted list item \"Twig\" right after the \"Npm 9.5.1\" item\n# in the \"Pour le front :\" list (same list/style as its neighbours).\n\n$d = $word.ActiveDocument\n\n$count = $d.Paragraphs.Count\n$targetIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text -like \"*9.5.1*\") {\n        $targetIndex = $i\n        break\n    }\n}\n\nif ($targetIndex -eq -1) {\n    Write-Output \"Could not find the '9.5.1' paragraph to anchor the new item.\"\n} else {\n    $target = $d.Paragraphs.Item($targetIndex)\n    # InsertParagraphAfter creates a new paragraph that inherits the source\n    # paragraph's style + numbering (pStyle \"Paragraphedeliste\", numId 2).\n    $target.Range.InsertParagraphAfter()\n    $newPara = $d.Paragraphs.Item($targetIndex + 1)\n    $newPara.Range.Text = \"Twig\"\n}\n"}
